$d = $word.ActiveDocument

# Locate the "Making the views" Heading-1 paragraph; the new table must be
# inserted right before it (i.e. right after the empty paragraph that
# currently precedes it).
$rng = $d.Content
$found = $rng.Find.Execute("Making the views", $true, $false, $false, $false,
                            $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate insertion point ('Making the views' heading)."
}

# Collapse the found range to its start so the table is inserted right
# before the heading, without replacing any existing text.
$insertionPoint = $d.Range($rng.Start, $rng.Start)

$xmlFragment = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:tbl>
<w:tblPr>
<w:tblStyle w:val="Tabelraster"/>
<w:tblW w:w="0" w:type="auto"/>
<w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/>
</w:tblPr>
<w:tblGrid>
<w:gridCol w:w="9062"/>
</w:tblGrid>
<w:tr>
<w:tc>
<w:tcPr>
<w:tcW w:w="9062" w:type="dxa"/>
</w:tcPr>
<w:p>
<w:pPr>
<w:rPr>
<w:lang w:val="en-GB"/>
</w:rPr>
</w:pPr>
<w:r>
<w:rPr>
<w:lang w:val="en-GB"/>
</w:rPr>
<w:t xml:space="preserve">Opening the admin-page of a </w:t>
</w:r>
<w:r>
<w:rPr>
<w:color w:val="EE0000"/>
<w:lang w:val="en-GB"/>
</w:rPr>
<w:t>SITE</w:t>
</w:r>
<w:r>
<w:rPr>
<w:lang w:val="en-GB"/>
</w:rPr>
<w:t>.</w:t>
</w:r>
</w:p>
</w:tc>
</w:tr>
<w:tr>
<w:tc>
<w:tcPr>
<w:tcW w:w="9062" w:type="dxa"/>
<w:shd w:val="clear" w:color="auto" w:fill="F2F2F2" w:themeFill="background1" w:themeFillShade="F2"/>
</w:tcPr>
<w:p>
<w:pPr>
<w:rPr>
<w:lang w:val="en-GB"/>
</w:rPr>
</w:pPr>
<w:r>
<w:rPr>
<w:lang w:val="en-GB"/>
</w:rPr>
<w:t>In your browser type</w:t>
</w:r>
<w:r>
<w:rPr>
<w:lang w:val="en-GB"/>
</w:rPr>
<w:t xml:space="preserve">: </w:t>
</w:r>
<w:r>
<w:rPr>
<w:color w:val="EE0000"/>
<w:lang w:val="en-GB"/>
</w:rPr>
<w:t>SITE-URL</w:t>
</w:r>
<w:r>
<w:rPr>
<w:lang w:val="en-GB"/>
</w:rPr>
<w:t>/admin</w:t>
</w:r>
</w:p>
</w:tc>
</w:tr>
</w:tbl>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$null = $insertionPoint.InsertXML($xmlFragment)
